$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = 1231234
$ws.Range("D2").Value = "johndoe@gmail.com"
$ws.Range("E2").Value = 123123123
$ws.Range("F2").Value = "Test Address"

# --- Row 3 ---
$ws.Range("A3").Value = 123123
$ws.Range("B3").Value = "John"
$ws.Range("C3").Value = "Doe"
$ws.Range("D3").Value = "johndoe@gmail.com"
$ws.Range("E3").Value = 123123123
$ws.Range("F3").Value = "Test Address"

# --- Row 4: cleared out entirely ---
$ws.Range("A4:F4").ClearContents()

# --- Hyperlinks on the two Email cells ---
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:johndoe@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:johndoe@gmail.com")

# --- Column widths (characters; target XML widths are 12.5/14.875/13/28.75/19.25/17.625) ---
$ws.Columns.Item(1).ColumnWidth = 11.666666666666666
$ws.Columns.Item(2).ColumnWidth = 14.0
$ws.Columns.Item(3).ColumnWidth = 12.166666666666666
$ws.Columns.Item(4).ColumnWidth = 28.0
$ws.Columns.Item(5).ColumnWidth = 18.5
$ws.Columns.Item(6).ColumnWidth = 16.833333333333332

# --- Selection ---
$ws.Range("E6").Select() | Out-Null
